$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete row 6 (data reduced from 5 data rows to 4; dimension becomes A1:AH5)
$ws.Rows.Item(6).Delete()

# 2) Narrow specific columns from width 8 to width 7 (M=13, Q=17, AA=27, AB=28, AC=29)
# Excel ColumnWidth differs from the raw OOXML <col width> by a constant offset (~0.83
# for this workbooks default font), determined empirically against the unmodified columns.
$ws.Columns.Item(13).ColumnWidth = 7 - 0.83
$ws.Columns.Item(17).ColumnWidth = 7 - 0.83
$ws.Columns.Item(27).ColumnWidth = 7 - 0.83
$ws.Columns.Item(28).ColumnWidth = 7 - 0.83
$ws.Columns.Item(29).ColumnWidth = 7 - 0.83

# 3) Replace the data in rows 2-5 (all columns A:AH) with the new dataset values
# Row 2
$ws.Cells.Item(2, 1).Value = 45157.50694444445
$ws.Cells.Item(2, 2).Value = 10.232
$ws.Cells.Item(2, 3).Value = 6.988
$ws.Cells.Item(2, 4).Value = 3.382
$ws.Cells.Item(2, 5).Value = 22.493
$ws.Cells.Item(2, 6).Value = 16.551
$ws.Cells.Item(2, 7).Value = 7.711
$ws.Cells.Item(2, 8).Value = 22.993
$ws.Cells.Item(2, 9).Value = 12.747
$ws.Cells.Item(2, 10).Value = 5.033
$ws.Cells.Item(2, 11).Value = 6.872
$ws.Cells.Item(2, 12).Value = 8.682
$ws.Cells.Item(2, 13).Value = 9.968
$ws.Cells.Item(2, 14).Value = 2.214
$ws.Cells.Item(2, 15).Value = 8.271000000000001
$ws.Cells.Item(2, 16).Value = 11.091
$ws.Cells.Item(2, 17).Value = 7.662
$ws.Cells.Item(2, 18).Value = 2.404
$ws.Cells.Item(2, 19).Value = 1.071
$ws.Cells.Item(2, 20).Value = 118.497
$ws.Cells.Item(2, 21).Value = 22.775
$ws.Cells.Item(2, 22).Value = 7.635
$ws.Cells.Item(2, 23).Value = 14.429
$ws.Cells.Item(2, 24).Value = 7.663
$ws.Cells.Item(2, 25).Value = 2.148
$ws.Cells.Item(2, 26).Value = 12.965
$ws.Cells.Item(2, 27).Value = 6.744
$ws.Cells.Item(2, 28).Value = 6.372
$ws.Cells.Item(2, 29).Value = 7.242
$ws.Cells.Item(2, 30).Value = 9.571
$ws.Cells.Item(2, 31).Value = 2.474
$ws.Cells.Item(2, 32).Value = 20.428
$ws.Cells.Item(2, 33).Value = 3.858
$ws.Cells.Item(2, 34).Value = 9.544

# Row 3
$ws.Cells.Item(3, 1).Value = 45157.51388888889
$ws.Cells.Item(3, 2).Value = 14.189
$ws.Cells.Item(3, 3).Value = 10.401
$ws.Cells.Item(3, 4).Value = 1.684
$ws.Cells.Item(3, 5).Value = 31.258
$ws.Cells.Item(3, 6).Value = 24.718
$ws.Cells.Item(3, 7).Value = 10.979
$ws.Cells.Item(3, 8).Value = 41.511
$ws.Cells.Item(3, 9).Value = 17.405
$ws.Cells.Item(3, 10).Value = 7.615
$ws.Cells.Item(3, 11).Value = 10.755
$ws.Cells.Item(3, 12).Value = 12.432
$ws.Cells.Item(3, 13).Value = 13.53
$ws.Cells.Item(3, 14).Value = 3.359
$ws.Cells.Item(3, 15).Value = 11.279
$ws.Cells.Item(3, 16).Value = 15.772
$ws.Cells.Item(3, 17).Value = 9.885
$ws.Cells.Item(3, 18).Value = 1.232
$ws.Cells.Item(3, 19).Value = 0.7
$ws.Cells.Item(3, 20).Value = 164.346
$ws.Cells.Item(3, 21).Value = 31.51
$ws.Cells.Item(3, 22).Value = 10.411
$ws.Cells.Item(3, 23).Value = 20.791
$ws.Cells.Item(3, 24).Value = 10.925
$ws.Cells.Item(3, 25).Value = 1.889
$ws.Cells.Item(3, 26).Value = 21.169
$ws.Cells.Item(3, 27).Value = 9.196
$ws.Cells.Item(3, 28).Value = 8.345000000000001
$ws.Cells.Item(3, 29).Value = 9.747
$ws.Cells.Item(3, 30).Value = 13.191
$ws.Cells.Item(3, 31).Value = 1.101
$ws.Cells.Item(3, 32).Value = 37.885
$ws.Cells.Item(3, 33).Value = 5.654
$ws.Cells.Item(3, 34).Value = 13.016

# Row 4
$ws.Cells.Item(4, 1).Value = 45157.52083333334
$ws.Cells.Item(4, 2).Value = 6.559
$ws.Cells.Item(4, 3).Value = 4.785
$ws.Cells.Item(4, 4).Value = 0.974
$ws.Cells.Item(4, 5).Value = 14.592
$ws.Cells.Item(4, 6).Value = 11.303
$ws.Cells.Item(4, 7).Value = 5.028
$ws.Cells.Item(4, 8).Value = 23.508
$ws.Cells.Item(4, 9).Value = 8.1
$ws.Cells.Item(4, 10).Value = 3.582
$ws.Cells.Item(4, 11).Value = 4.796
$ws.Cells.Item(4, 12).Value = 5.802
$ws.Cells.Item(4, 13).Value = 6.383
$ws.Cells.Item(4, 14).Value = 1.506
$ws.Cells.Item(4, 15).Value = 5.263
$ws.Cells.Item(4, 16).Value = 7.318
$ws.Cells.Item(4, 17).Value = 4.749
$ws.Cells.Item(4, 18).Value = 0.8080000000000001
$ws.Cells.Item(4, 19).Value = 0.349
$ws.Cells.Item(4, 20).Value = 72.81
$ws.Cells.Item(4, 21).Value = 14.87
$ws.Cells.Item(4, 22).Value = 4.858
$ws.Cells.Item(4, 23).Value = 9.702
$ws.Cells.Item(4, 24).Value = 5.055
$ws.Cells.Item(4, 25).Value = 0.944
$ws.Cells.Item(4, 26).Value = 11.452
$ws.Cells.Item(4, 27).Value = 4.291
$ws.Cells.Item(4, 28).Value = 3.962
$ws.Cells.Item(4, 29).Value = 4.618
$ws.Cells.Item(4, 30).Value = 6.156
$ws.Cells.Item(4, 31).Value = 0.708
$ws.Cells.Item(4, 32).Value = 21.669
$ws.Cells.Item(4, 33).Value = 2.572
$ws.Cells.Item(4, 34).Value = 6.075

# Row 5
$ws.Cells.Item(5, 1).Value = 45157.52777777778
$ws.Cells.Item(5, 2).Value = 3.71
$ws.Cells.Item(5, 3).Value = 2.69
$ws.Cells.Item(5, 4).Value = 0.67
$ws.Cells.Item(5, 5).Value = 8.34
$ws.Cells.Item(5, 6).Value = 6.32
$ws.Cells.Item(5, 7).Value = 2.81
$ws.Cells.Item(5, 8).Value = 14.15
$ws.Cells.Item(5, 9).Value = 4.61
$ws.Cells.Item(5, 10).Value = 2.07
$ws.Cells.Item(5, 11).Value = 2.61
$ws.Cells.Item(5, 12).Value = 3.32
$ws.Cells.Item(5, 13).Value = 3.69
$ws.Cells.Item(5, 14).Value = 0.82
$ws.Cells.Item(5, 15).Value = 3.01
$ws.Cells.Item(5, 16).Value = 4.14
$ws.Cells.Item(5, 17).Value = 2.79
$ws.Cells.Item(5, 18).Value = 0.61
$ws.Cells.Item(5, 19).Value = 0.21
$ws.Cells.Item(5, 20).Value = 38.47
$ws.Cells.Item(5, 21).Value = 8.529999999999999
$ws.Cells.Item(5, 22).Value = 2.78
$ws.Cells.Item(5, 23).Value = 5.5
$ws.Cells.Item(5, 24).Value = 2.86
$ws.Cells.Item(5, 25).Value = 0.58
$ws.Cells.Item(5, 26).Value = 6.81
$ws.Cells.Item(5, 27).Value = 2.45
$ws.Cells.Item(5, 28).Value = 2.3
$ws.Cells.Item(5, 29).Value = 2.68
$ws.Cells.Item(5, 30).Value = 3.52
$ws.Cells.Item(5, 31).Value = 0.52
$ws.Cells.Item(5, 32).Value = 13.01
$ws.Cells.Item(5, 33).Value = 1.43
$ws.Cells.Item(5, 34).Value = 3.47

